# Apply the "adding correlations, fixed the discretized attempts bug" edit:
#  - Add a new "Variable Definitions" worksheet at the end of the workbook
#  - Populate it with a small variable glossary table
#  - Tidy up the HurdleRates sheet (selection + explicit portrait page setup)
#  - Leave the newly added sheet as the active / selected tab

$wb = $excel.ActiveWorkbook

# --- HurdleRates sheet tweaks -------------------------------------------------
$wsHurdle = $wb.Worksheets.Item("HurdleRates")
$wsHurdle.PageSetup.Orientation = 1          # xlPortrait -> writes <pageSetup orientation="portrait"/>
$wsHurdle.Range("A15").Select() | Out-Null   # leave selection parked at A15 (no longer the active tab)

# --- New "Variable Definitions" sheet, appended after the last sheet ---------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsVars = $wb.Worksheets.Add($null, $lastSheet)
$wsVars.Name = "Variable Definitions"

# Fill the table in the same order the strings were originally authored in,
# so the shared-string table lines up with the source workbook:
#   header row -> variable symbols down column A -> descriptions (bottom-up)
#   -> final row appended last.
$wsVars.Range("A1").Value = "Variable Name"
$wsVars.Range("B1").Value = "Description"

$wsVars.Range("A2").Value = '$/theta$'
$wsVars.Range("A3").Value = '$n$'
$wsVars.Range("A4").Value = '$p$'

$wsVars.Range("B4").Value = "Per-attempt probability of success"
$wsVars.Range("B3").Value = "Number of attempts"
$wsVars.Range("B2").Value = "Target probability of success"

$wsVars.Range("A5").Value = '$X$'
$wsVars.Range("B5").Value = "Present value of pull size"

# Leave the selection/active tab on the new sheet, matching the saved view state
$wsVars.Range("E30").Select() | Out-Null
